$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1080.8667
$ws.Range("I98").Value = 1028.9474
$ws.Range("J98").Value = 1170.5454
$ws.Range("K98").Value = 1028.9474
$ws.Range("L98").Value = 1170.5454
$ws.Range("M98").Value = 469.0526
$ws.Range("N98").Value = -4166.5454

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1080.8667
$ws.Range("I122").Value = 1028.9474
$ws.Range("J122").Value = 1170.5454
$ws.Range("K122").Value = 3086.8422
$ws.Range("L122").Value = 3511.6362
$ws.Range("M122").Value = -636.8422
$ws.Range("N122").Value = -8411.636200000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 6103806.5
$ws.Range("I132").Value = 4560.794
$ws.Range("J132").Value = 35728716
$ws.Range("K132").Value = 13682.382
$ws.Range("L132").Value = 107186148
$ws.Range("M132").Value = -11152.382
$ws.Range("N132").Value = -107191208

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 5884373.5
$ws.Range("I137").Value = 1764.04
$ws.Range("J137").Value = 22224956
$ws.Range("K137").Value = 5292.12
$ws.Range("L137").Value = 66674868
$ws.Range("M137").Value = -2742.12
$ws.Range("N137").Value = -66679968

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4466331.5
$ws.Range("I138").Value = 1096.0244
$ws.Range("J138").Value = 16671308
$ws.Range("K138").Value = 3288.0732
$ws.Range("L138").Value = 50013924
$ws.Range("M138").Value = 1851.9268
$ws.Range("N138").Value = -50024204

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 43750
$ws.Range("J92").Value = 43750
$ws.Range("L92").Value = 43750
$ws.Range("N92").Value = -48742

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2696.183
$ws.Range("I134").Value = 2075.9482
$ws.Range("J134").Value = 5463.385
$ws.Range("K134").Value = 6227.844599999999
$ws.Range("L134").Value = 16390.155
$ws.Range("M134").Value = -3692.844599999999
$ws.Range("N134").Value = -21460.155

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H93").Value = 10062.8
$ws.Range("I93").Value = 7686.6665
$ws.Range("J93").Value = 31448
$ws.Range("K93").Value = 7686.6665
$ws.Range("L93").Value = 31448
$ws.Range("M93").Value = -5814.6665
$ws.Range("N93").Value = -35192

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 477.70834
$ws.Range("I107").Value = 386.16666
$ws.Range("J107").Value = 569.25
$ws.Range("K107").Value = 386.16666
$ws.Range("L107").Value = 569.25
$ws.Range("M107").Value = 1533.83334
$ws.Range("N107").Value = -4409.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 43.909092
$ws.Range("I12").Value = 24.666666
$ws.Range("J12").Value = 51.125
$ws.Range("K12").Value = 73.99999800000001
$ws.Range("L12").Value = 153.375
$ws.Range("M12").Value = 99.00000199999999
$ws.Range("N12").Value = -499.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 757.33
$ws.Range("I131").Value = 430.57144
$ws.Range("J131").Value = 810.52325
$ws.Range("K131").Value = 1291.71432
$ws.Range("L131").Value = 2431.56975
$ws.Range("M131").Value = 3748.28568
$ws.Range("N131").Value = -12511.56975

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 3476
$ws.Range("J134").Value = 5955.5557
$ws.Range("L134").Value = 17866.6671
$ws.Range("N134").Value = -28006.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 879.61536
$ws.Range("I22").Value = 870.7143
$ws.Range("J22").Value = 890
$ws.Range("K22").Value = 870.7143
$ws.Range("L22").Value = 890
$ws.Range("M22").Value = -575.7143
$ws.Range("N22").Value = -1480

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 879.61536
$ws.Range("I27").Value = 870.7143
$ws.Range("J27").Value = 890
$ws.Range("K27").Value = 870.7143
$ws.Range("L27").Value = 890
$ws.Range("M27").Value = -763.7143
$ws.Range("N27").Value = -1104

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8541.583000000001
$ws.Range("I40").Value = 9562.375
$ws.Range("J40").Value = 6500
$ws.Range("K40").Value = 9562.375
$ws.Range("L40").Value = 6500
$ws.Range("M40").Value = -9426.375
$ws.Range("N40").Value = -6772

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1992
$ws.Range("I68").Value = 1466.6666
$ws.Range("J68").Value = 2063.6365
$ws.Range("K68").Value = 1466.6666
$ws.Range("L68").Value = 2063.6365
$ws.Range("M68").Value = -717.6666
$ws.Range("N68").Value = -3561.6365

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1992
$ws.Range("I71").Value = 1466.6666
$ws.Range("J71").Value = 2063.6365
$ws.Range("K71").Value = 7333.333000000001
$ws.Range("L71").Value = 10318.1825
$ws.Range("M71").Value = -3589.333000000001
$ws.Range("N71").Value = -17806.1825

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 12202239
$ws.Range("I132").Value = 4251.579
$ws.Range("J132").Value = 22736864
$ws.Range("K132").Value = 12754.737
$ws.Range("L132").Value = 68210592
$ws.Range("M132").Value = -10224.737
$ws.Range("N132").Value = -68215652

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H134").Value = 14809
$ws.Range("J134").Value = 14809
$ws.Range("L134").Value = 14809
$ws.Range("N134").Value = -24949

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H137").Value = 39614.5
$ws.Range("J137").Value = 39614.5
$ws.Range("L137").Value = 39614.5
$ws.Range("N137").Value = -49814.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10435.412
$ws.Range("I62").Value = 5900
$ws.Range("J62").Value = 10718.875
$ws.Range("K62").Value = 5900
$ws.Range("L62").Value = 10718.875
$ws.Range("M62").Value = -5276
$ws.Range("N62").Value = -11966.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 10435.412
$ws.Range("I65").Value = 5900
$ws.Range("J65").Value = 10718.875
$ws.Range("K65").Value = 29500
$ws.Range("L65").Value = 53594.375
$ws.Range("M65").Value = -26380
$ws.Range("N65").Value = -59834.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1526.8334
$ws.Range("I107").Value = 2034.5714
$ws.Range("J107").Value = 816
$ws.Range("K107").Value = 6103.7142
$ws.Range("L107").Value = 2448
$ws.Range("M107").Value = -4183.7142
$ws.Range("N107").Value = -6288

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4293.6216
$ws.Range("I132").Value = 4685.032
$ws.Range("J132").Value = 2271.3333
$ws.Range("K132").Value = 14055.096
$ws.Range("L132").Value = 6813.999899999999
$ws.Range("M132").Value = -11525.096
$ws.Range("N132").Value = -11873.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 14055.625
$ws.Range("J135").Value = 14055.625
$ws.Range("L135").Value = 14055.625
$ws.Range("N135").Value = -24195.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 67440
$ws.Range("J137").Value = 67440
$ws.Range("L137").Value = 67440
$ws.Range("N137").Value = -77640

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H139").Value = 47500
$ws.Range("J139").Value = 47500
$ws.Range("L139").Value = 47500
$ws.Range("N139").Value = -57780

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 39990
$ws.Range("J141").Value = 39990
$ws.Range("L141").Value = 39990
$ws.Range("N141").Value = -50350
